$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new numeric-looking text must be forced to Text format
# so Excel does not auto-convert them to numbers (matching the original
# inline-string authoring of this worksheet).
$textFormatCells = @(
    "D5",
    "D6",
    "D14",
    "D20",
    "D21",
    "D24",
    "D25",
    "D26",
    "D29",
    "D30",
    "D32",
    "D33",
    "D35",
    "D36",
    "D38",
    "D42",
    "D46",
    "D48"
)
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated price / volume values cell by cell.
$ws.Range("D2").Value = "61.986.01"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "2.907.92"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "568.15"
$ws.Range("E5").Value = "  -3.34%  "
$ws.Range("D6").Value = "143.77"
$ws.Range("E6").Value = "  -1.47%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "2.904.97"
$ws.Range("E8").Value = "  -0.30%  "
$ws.Range("E9").Value = "  -1.12%  "
$ws.Range("E10").Value = "  -1.12%  "
$ws.Range("E11").Value = "  -1.91%  "
$ws.Range("E12").Value = "  -1.51%  "
$ws.Range("E13").Value = "  -0.60%  "
$ws.Range("D14").Value = "32.59"
$ws.Range("E14").Value = "  +0.35%  "
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").Value = "3.390.52"
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("D17").Value = "61.955.35"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").Value = "2.908.44"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("E19").Value = "  -1.53%  "
$ws.Range("D20").Value = "430.94"
$ws.Range("E20").Value = "  -0.85%  "
$ws.Range("D21").Value = "13.05"
$ws.Range("E21").Value = "  -2.93%  "
$ws.Range("E22").Value = "  -1.19%  "
$ws.Range("E23").Value = "  -1.12%  "
$ws.Range("D24").Value = "78.88"
$ws.Range("E24").Value = "  -2.46%  "
$ws.Range("D25").Value = "11.98"
$ws.Range("E25").Value = "  +0.74%  "
$ws.Range("D26").Value = "10.26"
$ws.Range("E26").Value = "  -5.54%  "
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("E28").Value = "  -3.33%  "
$ws.Range("D29").Value = "0.0000113"
$ws.Range("E29").Value = "  +8.06%  "
$ws.Range("D30").Value = "6.93"
$ws.Range("E30").Value = "  -5.46%  "
$ws.Range("E31").Value = "  -2.49%  "
$ws.Range("D32").Value = "2.00"
$ws.Range("E32").Value = "  -4.97%  "
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("E34").Value = "  -2.87%  "
$ws.Range("D35").Value = "25.67"
$ws.Range("E35").Value = "  -1.32%  "
$ws.Range("D36").Value = "0.954"
$ws.Range("E37").Value = "  -2.41%  "
$ws.Range("D38").Value = "48.82"
$ws.Range("E38").Value = "  -0.84%  "
$ws.Range("E39").Value = "  -6.15%  "
$ws.Range("E40").Value = "  -4.84%  "
$ws.Range("E41").Value = "  -0.72%  "
$ws.Range("D42").Value = "41.18"
$ws.Range("E42").Value = "  +5.57%  "
$ws.Range("E43").Value = "  -2.53%  "
$ws.Range("E44").Value = "  -2.58%  "
$ws.Range("D45").Value = "2.724.04"
$ws.Range("E45").Value = "  +0.91%  "
$ws.Range("D46").Value = "133.10"
$ws.Range("E46").Value = "  -1.22%  "
$ws.Range("E47").Value = "  -0.56%  "
$ws.Range("D48").Value = "349.42"
$ws.Range("E48").Value = "  +0.82%  "
$ws.Range("E50").Value = "  -1.12%  "
$ws.Range("E51").Value = "  +10.39%  "
